$d = $word.ActiveDocument

# Locate the paragraph that ends the "13)" answer ("The variables of JavaScript ...")
# so we can insert the new "15) ..." Q&A block right after it.
$rng = $d.Content
$found = $rng.Find.Execute(
    "The variables of JavaScript represent the arguments that are passed to a function.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Anchor paragraph not found"
}

$anchor = $rng.Paragraphs(1)

# --- Paragraph 1: "15) What's The Difference Between Undeclared & Undefined Variables?"
$anchor.Range.InsertParagraphAfter() | Out-Null
$pHeading = $anchor.Next()
$pHeading.Range.Text = "15) What's The Difference Between Undeclared & Undefined Variables?"

# Split "15)" / " " / "What's The Difference..." into separate runs, matching the
# original document's convention for numbered question headings.
$headingStart = $pHeading.Range.Start
$spaceRng = $d.Range($headingStart + 3, $headingStart + 4)
$spaceRng.Bold = 1
$spaceRng.Bold = 0

# --- Paragraph 2: undeclared-variable explanation
$pHeading.Range.InsertParagraphAfter() | Out-Null
$pUndeclared = $pHeading.Next()
$pUndeclared.Range.Text = "An undeclared variable has not been declared anywhere in the code, so said variable does not exist. If you try to read an undeclared variable, JavaScript throws an error."

# --- Paragraph 3: blank separator line
$pUndeclared.Range.InsertParagraphAfter() | Out-Null
$pBlank = $pUndeclared.Next()

# --- Paragraph 4: undefined-variable explanation
$pBlank.Range.InsertParagraphAfter() | Out-Null
$pUndefined = $pBlank.Next()
$pUndefined.Range.Text = "An undefined variable has been declared in the program, but no value has been assigned. This means the variable exists, but its value is yet to be defined."

Write-Output "Inserted Q15 block (undeclared vs undefined variables)."
